# Rename all 50 worksheets to new summ<number> identifiers, matching their
# existing position/order in the workbook (sheet index -> new name).
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ54712533",
    "summ54897501",
    "summ55072527",
    "summ55237460",
    "summ55397355",
    "summ55557804",
    "summ55715846",
    "summ55979836",
    "summ56142500",
    "summ56304194",
    "summ56464470",
    "summ56625587",
    "summ56786727",
    "summ56945568",
    "summ57103198",
    "summ57262185",
    "summ57422015",
    "summ57581491",
    "summ57738777",
    "summ57907832",
    "summ58091529",
    "summ58316791",
    "summ58557841",
    "summ58798182",
    "summ59015936",
    "summ59281600",
    "summ59535899",
    "summ59757300",
    "summ59987197",
    "summ00237289",
    "summ00447815",
    "summ00672968",
    "summ00920048",
    "summ01197301",
    "summ01448722",
    "summ01690276",
    "summ01946580",
    "summ02202052",
    "summ02533852",
    "summ02759857",
    "summ02992693",
    "summ03209343",
    "summ03420953",
    "summ03644051",
    "summ03868861",
    "summ04092660",
    "summ04341369",
    "summ04567791",
    "summ04801714",
    "summ05014488"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
